$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.230414390563965
$ws.Range("B1").Value = 3.267342567443848
$ws.Range("C1").Value = 4.606362342834473
$ws.Range("D1").Value = 2.185204267501831
$ws.Range("E1").Value = 1.526758551597595
